# Updates the cached "datetimeFigureOut" date placeholders (5/21/2021 -> 6/4/2021)
# across the slide layouts + handout master, and resizes/repositions the
# "Picture Placeholder 2" shape on the "Picture with Caption" layout
# (slide layout #9) to match the authored change.

$p = $ppt.ActivePresentation

$oldDate = "5/21/2021"
$newDate = "6/4/2021"

# --- 1. Refresh the cached date text on every slide layout that carries one ---
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $cl.Shapes.Count; $si++) {
        $sh = $cl.Shapes.Item($si)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Refresh the cached date text on the handout master ---
$handout = $p.HandoutMaster
for ($si = 1; $si -le $handout.Shapes.Count; $si++) {
    $sh = $handout.Shapes.Item($si)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 3. Resize/reposition the picture placeholder on "Picture with Caption" ---
$pictureLayout = $null
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    if ($cl.Name -eq "Picture with Caption") {
        $pictureLayout = $cl
    }
}

for ($si = 1; $si -le $pictureLayout.Shapes.Count; $si++) {
    $sh = $pictureLayout.Shapes.Item($si)
    if ($sh.Name -eq "Picture Placeholder 2") {
        $sh.Left = 19.18748031496063
        $sh.Top = 99.44165354330708
        $sh.Width = 523.2122047244095
        $sh.Height = 383.75
    }
}
